$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string country name swaps (rank changes) ---
$ws.Range("A78").Value = "Bosnia y Herzegovina"
$ws.Range("A79").Value = "Dinamarca"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Update timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 14:06"

# --- Update numeric case data ---
# Row 4
$ws.Range("B4").Value = 5416040
$ws.Range("C4").Value = 374
$ws.Range("E4").Value = 2401976

# Row 6
$ws.Range("B6").Value = 2465662
$ws.Range("C6").Value = 6049
$ws.Range("D6").Value = 1753719
$ws.Range("E6").Value = 663755
$ws.Range("G6").Value = 44
$ws.Range("H6").Value = 48188

# Row 38
$ws.Range("B38").Value = 82743
$ws.Range("C38").Value = 212
$ws.Range("D38").Value = 77427
$ws.Range("E38").Value = 4759
$ws.Range("G38").Value = 6
$ws.Range("H38").Value = 557

# Row 77
$ws.Range("B77").Value = 15834
$ws.Range("C77").Value = 343
$ws.Range("D77").Value = 9382
$ws.Range("E77").Value = 6346

# Row 78
$ws.Range("B78").Value = 15535
$ws.Range("C78").Value = 351
$ws.Range("D78").Value = 9344
$ws.Range("E78").Value = 5722
$ws.Range("G78").Value = 11
$ws.Range("H78").Value = 469

# Row 79
$ws.Range("B79").Value = 15379
$ws.Range("C79").Value = 165
$ws.Range("D79").Value = 13216
$ws.Range("E79").Value = 1542
$ws.Range("H79").Value = 621

# Row 82
$ws.Range("B82").Value = 13643
$ws.Range("C82").Value = 121
$ws.Range("D82").Value = 12011
$ws.Range("E82").Value = 1468
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 164

# Row 84
$ws.Range("B84").Value = 12162
$ws.Range("C84").Value = 47
$ws.Range("D84").Value = 6325
$ws.Range("E84").Value = 5044
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 793

# Row 86
$ws.Range("E86").Value = 733
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 261

# Row 128
$ws.Range("D128").Value = 2027
$ws.Range("E128").Value = 213

# Row 134
$ws.Range("B134").Value = 1983
$ws.Range("C134").Value = 7
$ws.Range("D134").Value = 1861
$ws.Range("E134").Value = 112

# Row 158
$ws.Range("B158").Value = 929
$ws.Range("C158").Value = 18
$ws.Range("D158").Value = 437
$ws.Range("E158").Value = 471

# Row 206
$ws.Range("B206").Value = 22
$ws.Range("C206").Value = 2
$ws.Range("E206").Value = 3

# Row 213
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

